# Chapter 4 ("Protein folding and artificial intelligence") contains a
# sub-heading "NP-completeness and NP-hardness" that, unlike its sibling
# heading "Definition and mechanism", is not yet colored red. Make it red
# (matching the sibling), i.e. add <w:color w:val="FF0000"/> to the run
# properties of the paragraph mark and of both text runs in that paragraph.

$d = $word.ActiveDocument

$targetText = "NP-completeness and NP-hardness"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    # Strip the trailing paragraph mark / cell mark characters so the
    # comparison is exact regardless of how the mark is represented.
    $t = $t.Replace([char]13, "").Replace([char]7, "")

    if ($t -eq $targetText) {
        # Setting Font.Color on the paragraph's own Range (which includes
        # the paragraph mark) updates the rPr of the paragraph mark as well
        # as the rPr of every run contained in the paragraph.
        $p.Range.Font.Color = 255
    }
}
